# Applies the two changes captured by the commit:
#   1. Slide 5's table switches from the custom "Table_0" style to the
#      built-in "No Style, No Grid" table style.
#   2. The deck's (editable) theme color scheme is swapped from the
#      "Integral" design's "Red Violet" palette to the stock "Office"
#      palette (i.e. the presentation's look becomes the default Office
#      Theme colors).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{8C073611-9E7F-4288-9D49-18983821D71F}")

# --- 2. Theme color scheme -> stock "Office" palette ----------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Index order matches <a:clrScheme>: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
